$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.006.01"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.051.98"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.35%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.052.06"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.153"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.13%  "
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.550.88"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.14"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.012.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.051.52"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.48"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.50"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.93"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.65"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.53"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.71%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.65"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.68"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0814"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("E37").Value = "  -3.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.22"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.25"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.41"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "428.65"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.286"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.840.47"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.97%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.09"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -1.52%  "
